$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("product backlog")
Write-Host $ws1.Name
